# Re-arrange the "No Stock" item list: within each brand group the
# individual item rows (Item Name + UOM) are reordered to match the
# author's new layout. BSL NO / BRAND / ISL NO / Total Ordered /
# Estimated Sales stay attached to their original row position; only the
# Item Name (col D) and UOM (col E) values move, paired together so the
# UOM always follows its item.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 27

# Snapshot the current Item Name -> UOM pairing (order doesn't matter,
# we look values up by item name so nothing is lost or mismatched).
$itemToUom = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $item = $ws.Cells.Item($r, 4).Value()
    $uom = $ws.Cells.Item($r, 5).Value()
    $itemToUom[$item] = $uom
}

# Target order of Item Name for rows 2..27 (brand grouping is unchanged,
# only the order of items within each brand group changes).
$newItemOrder = @(
    "Desodin 60ml Syrup",
    "Dinafex 120mg Tablet",
    "Dinafex 180mg Tablet",
    "Dinafex 60mg Tablet",
    "Dorenta 50mg Tablet",
    "Etorix 90mg Tablet",
    "Etorix 120mg Tablet",
    "Etorix 60mg Tablet - 40's",
    "Fenobac 100ml Syrup",
    "Flucloxin 500mg Capsule - 36's",
    "Flucloxin 500mg Capsule",
    "Geminox 320mg Tablet - 8's",
    "Ketonic 10mg Tablet",
    "Ketonic 30mg Injection",
    "Ketonic 30mg IM/IV Injection - 4's",
    "Kynol D 25mg Tablet",
    "Kynol TR 200mg Capsule",
    "Kynol TR 100mg Capsule",
    "Naprox Plus 500mg Tablet - 30's",
    "Oradin Plus Tablet - 40's",
    "Osticare Tablet 24's",
    "Rupaday Oral Solution 60ml",
    "Zithrox 15ml Suspension",
    "Zithrox 250mg Tablet - 6's",
    "Zithrox 500mg Tablet",
    "Zithrox 30ml Dry Suspension"
)

for ($i = 0; $i -lt $newItemOrder.Length; $i++) {
    $r = $firstRow + $i
    $item = $newItemOrder[$i]
    $ws.Cells.Item($r, 4).Value = $item
    $ws.Cells.Item($r, 5).Value = $itemToUom[$item]
}
